$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns are stored as text in the source data (e.g. "51.674.90",
# "  -0.19%  "); force Text format so Excel does not reinterpret them as numbers
# and strip formatting like trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '51.674.90'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '2.818.98'
$ws.Range("E3").Value = '  +2.24%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = '351.60'
$ws.Range("E5").Value = '  +5.78%  '

$ws.Range("D6").Value = '113.45'
$ws.Range("E6").Value = '  -2.28%  '

$ws.Range("D7").Value = '0.552'
$ws.Range("E7").Value = '  +2.53%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = '0.598'
$ws.Range("E9").Value = '  +3.96%  '

$ws.Range("D10").Value = '41.51'
$ws.Range("E10").Value = '  -0.51%  '

$ws.Range("D11").Value = '0.0848'
$ws.Range("E11").Value = '  -1.30%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.131'
$ws.Range("E12").Value = '  +1.42%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '19.95'
$ws.Range("E13").Value = '  -1.00%  '

$ws.Range("D14").Value = '7.70'
$ws.Range("E14").Value = '  +0.92%  '

$ws.Range("D15").Value = '3.276.75'
$ws.Range("E15").Value = '  +2.46%  '

$ws.Range("D16").Value = '2.847.26'
$ws.Range("E16").Value = '  +2.95%  '

$ws.Range("D17").Value = '0.892'
$ws.Range("E17").Value = '  +0.73%  '

$ws.Range("D18").Value = '51.743.82'
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").Value = '7.34'
$ws.Range("E19").Value = '  +6.83%  '

$ws.Range("D20").Value = '3.13'
$ws.Range("E20").Value = '  -2.72%  '

$ws.Range("D21").Value = '13.46'
$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").Value = '0.0₃0993'
$ws.Range("E22").Value = '  +1.52%  '

$ws.Range("D23").Value = '269.48'
$ws.Range("E23").Value = '  -3.09%  '

$ws.Range("D24").Value = '69.63'
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").Value = '2.78'
$ws.Range("E25").Value = '  +4.23%  '

$ws.Range("D26").Value = '26.59'
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").Value = '10.27'
$ws.Range("E28").Value = '  +0.92%  '

$ws.Range("E29").Value = '  +1.25%  '

$ws.Range("E30").Value = '  -1.15%  '

$ws.Range("B31").Value = 'OKB'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D31").Value = '50.64'
$ws.Range("E31").Value = '  +1.33%  '

$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '33.58'
$ws.Range("E32").Value = '  -4.04%  '

$ws.Range("B33").Value = 'VeChain'
$ws.Range("C33").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D33").Value = '0.0449'
$ws.Range("E33").Value = '  +27.83%  '

$ws.Range("D34").Value = '5.80'
$ws.Range("E34").Value = '  +4.32%  '

$ws.Range("D35").Value = '0.0826'
$ws.Range("E35").Value = '  +0.45%  '

$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("D37").Value = '2.07'
$ws.Range("E37").Value = '  -0.57%  '

$ws.Range("D38").Value = '4.88'
$ws.Range("E38").Value = '  -2.89%  '

$ws.Range("D39").Value = '3.21'
$ws.Range("E39").Value = '  -0.72%  '

$ws.Range("D40").Value = '18.03'
$ws.Range("E40").Value = '  -5.05%  '

$ws.Range("D41").Value = '23.74'
$ws.Range("E41").Value = '  +3.21%  '

$ws.Range("D42").Value = '2.55'
$ws.Range("E42").Value = '  +4.00%  '

$ws.Range("E43").Value = '  +0.80%  '

$ws.Range("D44").Value = '125.67'
$ws.Range("E44").Value = '  -1.10%  '

$ws.Range("E45").Value = '  +0.40%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.080.82'
$ws.Range("E46").Value = '  -0.29%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  +0.75%  '

$ws.Range("E48").Value = '  +3.71%  '

$ws.Range("D49").Value = '5.68'
$ws.Range("E49").Value = '  +2.58%  '

$ws.Range("D50").Value = '0.936'
$ws.Range("E50").Value = '  +6.81%  '

$ws.Range("D51").Value = '60.44'
$ws.Range("E51").Value = '  +0.79%  '
